$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text: "04-12-2025" -> "18-12-2025"
#    (Slide Master + all 11 Custom Layouts)
# ---------------------------------------------------------------------------
$newDate = "18-12-2025"

for ($j = 1; $j -le $p.SlideMaster.Shapes.Count; $j++) {
    $sh = $p.SlideMaster.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 title textbox ("TextBox 9", inside "Group 3"):
#    - reposition / resize so the longer title still fits
#    - retitle "CSK IPL ANALYSIS" -> "CSK IPL ANALYTICS"
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$title = $null
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $candidate = $slide.Shapes.Item($k)
    if ($candidate.Name -eq "Group 3") {
        for ($m = 1; $m -le $candidate.GroupItems.Count; $m++) {
            $inner = $candidate.GroupItems.Item($m)
            if ($inner.Name -eq "TextBox 9") {
                $title = $inner
            }
        }
    }
}

$title.Left = 68.8704337007874
$title.Top = 17.779646299212597
$title.Width = 375.74822897637796
$title.Height = 29.081299842519684

$title.TextFrame.TextRange.Text = "CSK IPL ANALYTICS"
